# Weekly update: insert a new "Fruta / hortaliza, semanal" observation at the
# top of the Mango / Vega Modelo de Temuco data block (row 199). Every
# existing record from row 199 down shifts one row lower, and a new row 241
# is created holding what used to be the last record (old row 240).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 199:240 down by one row, pushing the whole block (including the
# former last row) down to make room for the new record at row 199.
$ws.Rows.Item(199).Insert()

# Populate the freshly inserted row 199 with this week's new observation.
$ws.Cells.Item(199, 1).Value = 10
$ws.Cells.Item(199, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(199, 3).Value = "La Araucanía"
$ws.Cells.Item(199, 4).Value = 44508
$ws.Cells.Item(199, 5).Value = 9
$ws.Cells.Item(199, 6).Value = "Fruta"
$ws.Cells.Item(199, 7).Value = 100108
$ws.Cells.Item(199, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(199, 9).Value = 100108002
$ws.Cells.Item(199, 10).Value = "Mango"
$ws.Cells.Item(199, 11).Value = "Sin especificar"
$ws.Cells.Item(199, 12).Value = "Primera"
$ws.Cells.Item(199, 13).Value = 275
$ws.Cells.Item(199, 14).Value = 6000
$ws.Cells.Item(199, 15).Value = 7000
$ws.Cells.Item(199, 16).Value = 6545
$ws.Cells.Item(199, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(199, 18).Value = "Perú"
$ws.Cells.Item(199, 19).Value = 1636
$ws.Cells.Item(199, 20).Value = 4
